# The author opened Point.xlsx, appended a " (5đ)" price hint to the six
# item-name headers on the "CẮM TRẠI" sheet (Sheet2), then clicked on H3
# (selecting the merged H3:H4 header cell) before saving.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

$ws.Range("B3").Value = "Bàn (5đ)"
$ws.Range("C3").Value = "Ảnh Chúa (5đ)"
$ws.Range("D3").Value = "Kinh Thánh (5đ)"
$ws.Range("E3").Value = "Hoa (5đ)"
$ws.Range("F3").Value = "Bạt Ngồi (5đ)"
$ws.Range("G3").Value = "Cổng (5đ)"

$ws.Activate()
$ws.Range("H3:H4").Select()
